$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the now-obsolete DK*_ImportFrom_* rows (rows 36-41) and update
# the surviving rows in place with the aggregated/reordered tech -> modelTech
# mapping.
$ws.Rows("36:41").Delete()

$ws.Cells.Item(2, 1).Value = "BH_Biogas"
$ws.Cells.Item(2, 2).Value = "standard_H"
$ws.Cells.Item(3, 1).Value = "BP_Biogas"
$ws.Cells.Item(3, 2).Value = "BP"
$ws.Cells.Item(4, 1).Value = "IndustryH_Biogas"
$ws.Cells.Item(4, 2).Value = "standard_H"
$ws.Cells.Item(5, 1).Value = "BH_Biomass"
$ws.Cells.Item(5, 2).Value = "standard_H"
$ws.Cells.Item(6, 1).Value = "BP_Biomass"
$ws.Cells.Item(6, 2).Value = "BP"
$ws.Cells.Item(7, 1).Value = "BP_Coal"
$ws.Cells.Item(7, 2).Value = "BP"
$ws.Cells.Item(8, 1).Value = "BH_Natgas"
$ws.Cells.Item(8, 2).Value = "standard_H"
$ws.Cells.Item(9, 1).Value = "BP_Natgas"
$ws.Cells.Item(9, 2).Value = "BP"
$ws.Cells.Item(10, 1).Value = "IndustryH_Natgas"
$ws.Cells.Item(10, 2).Value = "standard_H"
$ws.Cells.Item(11, 1).Value = "BH_Oil"
$ws.Cells.Item(11, 2).Value = "standard_H"
$ws.Cells.Item(12, 1).Value = "BP_Oil"
$ws.Cells.Item(12, 2).Value = "BP"
$ws.Cells.Item(13, 1).Value = "IndustryH_Oil"
$ws.Cells.Item(13, 2).Value = "standard_H"
$ws.Cells.Item(14, 1).Value = "BH_Waste"
$ws.Cells.Item(14, 2).Value = "standard_H"
$ws.Cells.Item(15, 1).Value = "BP_Waste"
$ws.Cells.Item(15, 2).Value = "BP"
$ws.Cells.Item(16, 1).Value = "EP"
$ws.Cells.Item(16, 2).Value = "HP"
$ws.Cells.Item(17, 1).Value = "HPstandard"
$ws.Cells.Item(17, 2).Value = "HP"
$ws.Cells.Item(18, 1).Value = "HPsurplusheat"
$ws.Cells.Item(18, 2).Value = "HP"
$ws.Cells.Item(19, 1).Value = "IH"
$ws.Cells.Item(19, 2).Value = "HP"
$ws.Cells.Item(20, 1).Value = "IndustryH"
$ws.Cells.Item(20, 2).Value = "standard_H"
$ws.Cells.Item(21, 1).Value = "SH"
$ws.Cells.Item(21, 2).Value = "standard_H"
$ws.Cells.Item(22, 1).Value = "IndustryH_Biomass"
$ws.Cells.Item(22, 2).Value = "standard_H"
$ws.Cells.Item(23, 1).Value = "GT"
$ws.Cells.Item(23, 2).Value = "standard_H"
$ws.Cells.Item(24, 1).Value = "CD_Biogas"
$ws.Cells.Item(24, 2).Value = "standard_E"
$ws.Cells.Item(25, 1).Value = "IndustryE_Biogas"
$ws.Cells.Item(25, 2).Value = "standard_E"
$ws.Cells.Item(26, 1).Value = "IndustryE_Biomass"
$ws.Cells.Item(26, 2).Value = "standard_E"
$ws.Cells.Item(27, 1).Value = "CD_Coal"
$ws.Cells.Item(27, 2).Value = "standard_E"
$ws.Cells.Item(28, 1).Value = "CD_Natgas"
$ws.Cells.Item(28, 2).Value = "standard_E"
$ws.Cells.Item(29, 1).Value = "IndustryE_Natgas"
$ws.Cells.Item(29, 2).Value = "standard_E"
$ws.Cells.Item(30, 1).Value = "CD_Oil"
$ws.Cells.Item(30, 2).Value = "standard_E"
$ws.Cells.Item(31, 1).Value = "IndustryE_Oil"
$ws.Cells.Item(31, 2).Value = "standard_E"
$ws.Cells.Item(32, 1).Value = "PV"
$ws.Cells.Item(32, 2).Value = "standard_E"
$ws.Cells.Item(33, 1).Value = "ROR"
$ws.Cells.Item(33, 2).Value = "standard_E"
$ws.Cells.Item(34, 1).Value = "WL"
$ws.Cells.Item(34, 2).Value = "standard_E"
$ws.Cells.Item(35, 1).Value = "WS"
$ws.Cells.Item(35, 2).Value = "standard_E"
